{"js": "// Update page title, intro \"What we like\" / \"What we don't like\" bullet\n// points, and the closing SEO title/description paragraphs to match the\n// refreshed copy (\"Added many more features\").\n\nconst replacements = [\n  [\n    \"Play Cygnus Free Slot Game Online | ELK Studios\",\n    \"Play Cygnus - Free Online Slot Game\",\n  ],\n  [\n    \"A unique 'Avalanche' system with up to 262,144 ways to win\",\n    \"Unique 'Avalanche' system with up to 262,144 ways to win\",\n  ],\n  [\n    \"High-value symbol payouts, interesting multipliers\",\n    \"High value symbol payouts\",\n  ],\n  [\n    \"Fairly regular wins with substantial amounts with an RTP of 96.1%\",\n    \"Interesting multipliers\",\n  ],\n  [\n    \"Innovative and unique game design by reputable developer ELK Studios\",\n    \"Reputable developer with innovative game designs\",\n  ],\n  [\n    \"High volatility rating may not be suitable for players who prefer frequent wins\",\n    \"Infrequent wins due to high volatility\",\n  ],\n  [\n    \"Limited number of symbols and lack of bonus features may not be appealing to some players\",\n    \"Limited variety in symbol designs\",\n  ],\n  [\n    \"Read our review of Cygnus by ELK Studios and play this unique online slot game for free. Featuring up to 262,144 ways to win and a high volatility rating.\",\n    \"Discover the unique gameplay and high volatility of Cygnus. Play for free and win big!\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    // InsertLocation.replace overwrites only the matched range's text,\n    // keeping the surrounding runs/paragraph formatting (bold/italic, list\n    // style, etc.) untouched.\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update page title, intro \"What we like\" / \"What we don't like\" bullet\n# points, and the closing SEO title/description paragraphs to match the\n# refreshed copy (\"Added many more features\").\n\n$d = $word.ActiveDocument\n\n# Map of exact original paragraph text -> new paragraph text. Using a direct\n# Range.Text replacement (scoped to exclude the trailing paragraph mark)\n# keeps existing run formatting (bold/italic) intact and avoids the\n# smart-quote substitution that Find/Replace applies to straight quotes.\n$replacements = [ordered]@{\n    \"Play Cygnus Free Slot Game Online | ELK Studios\" = \"Play Cygnus - Free Online Slot Game\"\n    \"A unique 'Avalanche' system with up to 262,144 ways to win\" = \"Unique 'Avalanche' system with up to 262,144 ways to win\"\n    \"High-value symbol payouts, interesting multipliers\" = \"High value symbol payouts\"\n    \"Fairly regular wins with substantial amounts with an RTP of 96.1%\" = \"Interesting multipliers\"\n    \"Innovative and unique game design by reputable developer ELK Studios\" = \"Reputable developer with innovative game designs\"\n    \"High volatility rating may not be suitable for players who prefer frequent wins\" = \"Infrequent wins due to high volatility\"\n    \"Limited number of symbols and lack of bonus features may not be appealing to some players\" = \"Limited variety in symbol designs\"\n    \"Read our review of Cygnus by ELK Studios and play this unique online slot game for free. Featuring up to 262,144 ways to win and a high volatility rating.\" = \"Discover the unique gameplay and high volatility of Cygnus. Play for free and win big!\"\n}\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $full = $r.Text\n    if ($full.Length -gt 0) {\n        $plain = $full.Substring(0, $full.Length - 1)\n    } else {\n        $plain = $full\n    }\n    if ($replacements.Contains($plain)) {\n        $newText = $replacements[$plain]\n        $target = $d.Range($r.Start, $r.End - 1)\n        $target.Text = $newText\n    }\n}\n"}
